# Trade #48 closed at 2026-02-17 12:48:48 - unknown UNKNOWN +0.000%
#
# Updates the "Summary", "Strategy Status", "All Trades" and "MarketMaking"
# sheets to reflect the newly closed trade #48.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.16   # Current Capital
$summary.Range("B4").Value = 0.15      # Total P&L $
$summary.Range("B6").Value = 48        # Total Trades
$summary.Range("B7").Value = 21        # Winning Trades
$summary.Range("B9").Value = 43.75     # Win Rate %

# ---------------------------------------------------------------------------
# 2) Strategy Status sheet (MarketMaking row, row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.16     # Capital
$status.Range("D4").Value = 48         # Trades
$status.Range("E4").Value = 0.15       # P&L $
$status.Range("F4").Value = 0.16       # P&L %
$status.Range("G4").Value = 43.75      # Win Rate %

# ---------------------------------------------------------------------------
# 3) Helper to append the new trade record (row 49) to a trade-log sheet
# ---------------------------------------------------------------------------
function Add-Trade49Row($ws) {
    # Trade # (numeric)
    $ws.Cells.Item(49, 1).Value = 48

    # Date / Time are stored as plain text in this workbook, not Excel
    # dates/times, so force a text number format before assigning the
    # value to stop them being auto-converted into date/time serials.
    # ClearFormats() afterwards drops the temporary "@" number format
    # again so the cell is left with the default style, exactly like
    # the rest of the sheet.
    $ws.Cells.Item(49, 2).NumberFormat = "@"
    $ws.Cells.Item(49, 2).Value = "2026-02-17"
    $ws.Cells.Item(49, 2).ClearFormats()

    $ws.Cells.Item(49, 3).NumberFormat = "@"
    $ws.Cells.Item(49, 3).Value = "12:48:41"
    $ws.Cells.Item(49, 3).ClearFormats()

    $ws.Cells.Item(49, 4).Value = "MarketMaking"
    $ws.Cells.Item(49, 5).Value = "UP"
    $ws.Cells.Item(49, 6).Value = 0.94
    $ws.Cells.Item(49, 7).Value = 0.95
    $ws.Cells.Item(49, 8).Value = "CLOSED"
    $ws.Cells.Item(49, 9).Value = 1.0638
    $ws.Cells.Item(49, 10).Value = 0.01
    $ws.Cells.Item(49, 11).Value = 100.16
    $ws.Cells.Item(49, 12).Value = 0
    $ws.Cells.Item(49, 13).Value = 0
    $ws.Cells.Item(49, 14).Value = 0.6
    $ws.Cells.Item(49, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(49, 16).Value = "early_exit"
    $ws.Cells.Item(49, 17).Value = 0.11
}

# ---------------------------------------------------------------------------
# 4) All Trades sheet
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-Trade49Row $allTrades

# ---------------------------------------------------------------------------
# 5) MarketMaking sheet
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-Trade49Row $marketMaking
